$wb = $excel.ActiveWorkbook

$baseUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/95205c393d25374ad9cb12f6ac379006c040a924/e2e/"

# ---------------------------------------------------------------------------
# New file records being added by this handoff.
# ---------------------------------------------------------------------------
$file1 = "12f9654f-efa7-4e2e-a190-84c4eb579ed7.md"
$file2 = "63373375-e236-49d5-ad70-28fdb75af732.md"

$zhXlf1 = "12f9654f-efa7-4e2e-a190-84c4eb579ed7.1efdcc0cf4aa9ab79ef73c938af912cf49c0c6c4.zh-cn.xlf"
$zhXlf2 = "63373375-e236-49d5-ad70-28fdb75af732.4255cca06e51b170ad8ca4baee01702143772a43.zh-cn.xlf"
$deXlf1 = "12f9654f-efa7-4e2e-a190-84c4eb579ed7.1efdcc0cf4aa9ab79ef73c938af912cf49c0c6c4.de-de.xlf"
$deXlf2 = "63373375-e236-49d5-ad70-28fdb75af732.4255cca06e51b170ad8ca4baee01702143772a43.de-de.xlf"

$zhHandoffDt = "2016-10-18 10:35:01"
$deHandoffDt = "2016-10-18 10:35:25"
$status = "Ready for handoff"

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")
$lo = $ws.ListObjects.Item(1)
$lo.ListRows.Add() | Out-Null
$lo.ListRows.Add() | Out-Null

$ws.Cells.Item(4, 1).Value = $file1
$ws.Cells.Item(4, 3).Value = ".md"
$ws.Cells.Item(4, 4).Value = ""
$ws.Cells.Item(4, 5).Value = $status
$ws.Cells.Item(4, 6).Value = $status
$ws.Cells.Item(4, 7).Value = "2016-10-18 10:35:25"
$ws.Hyperlinks.Add($ws.Cells.Item(4, 2), ($baseUrl + $file1), [Type]::Missing, [Type]::Missing, ("e2e\" + $file1)) | Out-Null

$ws.Cells.Item(5, 1).Value = $file2
$ws.Cells.Item(5, 3).Value = ".md"
$ws.Cells.Item(5, 4).Value = ""
$ws.Cells.Item(5, 5).Value = $status
$ws.Cells.Item(5, 6).Value = $status
$ws.Cells.Item(5, 7).Value = "2016-10-18 10:35:25"
$ws.Hyperlinks.Add($ws.Cells.Item(5, 2), ($baseUrl + $file2), [Type]::Missing, [Type]::Missing, ("e2e\" + $file2)) | Out-Null

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")
$lo = $ws.ListObjects.Item(1)
$lo.ListRows.Add() | Out-Null
$lo.ListRows.Add() | Out-Null

$ws.Cells.Item(4, 2).Value = ".md"
$ws.Cells.Item(4, 3).Value = $status
$ws.Cells.Item(4, 4).Value = "e2e"
$ws.Cells.Item(4, 5).Value = "ht"
$ws.Cells.Item(4, 6).Value = "False"
$ws.Cells.Item(4, 7).Value = $zhXlf1
$ws.Cells.Item(4, 8).Value = $zhHandoffDt
$ws.Cells.Item(4, 9).Value = ""
$ws.Cells.Item(4, 10).Value = ""
$ws.Cells.Item(4, 11).Value = "0001-01-01 00:00:00"
$ws.Cells.Item(4, 12).Value = ""
$ws.Cells.Item(4, 13).Value = "True"
$ws.Cells.Item(4, 14).Value = ""
$ws.Cells.Item(4, 15).Value = "False"
$ws.Cells.Item(4, 16).Value = ""
$ws.Hyperlinks.Add($ws.Cells.Item(4, 1), ($baseUrl + $file1), [Type]::Missing, [Type]::Missing, $file1) | Out-Null

$ws.Cells.Item(5, 2).Value = ".md"
$ws.Cells.Item(5, 3).Value = $status
$ws.Cells.Item(5, 4).Value = "e2e"
$ws.Cells.Item(5, 5).Value = "ht"
$ws.Cells.Item(5, 6).Value = "False"
$ws.Cells.Item(5, 7).Value = $zhXlf2
$ws.Cells.Item(5, 8).Value = $zhHandoffDt
$ws.Cells.Item(5, 9).Value = ""
$ws.Cells.Item(5, 10).Value = ""
$ws.Cells.Item(5, 11).Value = "0001-01-01 00:00:00"
$ws.Cells.Item(5, 12).Value = ""
$ws.Cells.Item(5, 13).Value = "True"
$ws.Cells.Item(5, 14).Value = ""
$ws.Cells.Item(5, 15).Value = "False"
$ws.Cells.Item(5, 16).Value = ""
$ws.Hyperlinks.Add($ws.Cells.Item(5, 1), ($baseUrl + $file2), [Type]::Missing, [Type]::Missing, $file2) | Out-Null

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")
$lo = $ws.ListObjects.Item(1)
$lo.ListRows.Add() | Out-Null
$lo.ListRows.Add() | Out-Null

$ws.Cells.Item(4, 2).Value = ".md"
$ws.Cells.Item(4, 3).Value = $status
$ws.Cells.Item(4, 4).Value = "e2e"
$ws.Cells.Item(4, 5).Value = "ht"
$ws.Cells.Item(4, 6).Value = "False"
$ws.Cells.Item(4, 7).Value = $deXlf1
$ws.Cells.Item(4, 8).Value = $deHandoffDt
$ws.Cells.Item(4, 9).Value = ""
$ws.Cells.Item(4, 10).Value = ""
$ws.Cells.Item(4, 11).Value = "0001-01-01 00:00:00"
$ws.Cells.Item(4, 12).Value = ""
$ws.Cells.Item(4, 13).Value = "True"
$ws.Cells.Item(4, 14).Value = ""
$ws.Cells.Item(4, 15).Value = "False"
$ws.Cells.Item(4, 16).Value = ""
$ws.Hyperlinks.Add($ws.Cells.Item(4, 1), ($baseUrl + $file1), [Type]::Missing, [Type]::Missing, $file1) | Out-Null

$ws.Cells.Item(5, 2).Value = ".md"
$ws.Cells.Item(5, 3).Value = $status
$ws.Cells.Item(5, 4).Value = "e2e"
$ws.Cells.Item(5, 5).Value = "ht"
$ws.Cells.Item(5, 6).Value = "False"
$ws.Cells.Item(5, 7).Value = $deXlf2
$ws.Cells.Item(5, 8).Value = $deHandoffDt
$ws.Cells.Item(5, 9).Value = ""
$ws.Cells.Item(5, 10).Value = ""
$ws.Cells.Item(5, 11).Value = "0001-01-01 00:00:00"
$ws.Cells.Item(5, 12).Value = ""
$ws.Cells.Item(5, 13).Value = "True"
$ws.Cells.Item(5, 14).Value = ""
$ws.Cells.Item(5, 15).Value = "False"
$ws.Cells.Item(5, 16).Value = ""
$ws.Hyperlinks.Add($ws.Cells.Item(5, 1), ($baseUrl + $file2), [Type]::Missing, [Type]::Missing, $file2) | Out-Null
